# Updated cryptos list on Sun Nov  5 04:35:49 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures on the
# active worksheet to the latest coinranking.com snapshot, and corrects the
# ordering of the Gas / RenderToken rows (44 and 45).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.628.62'
$ws.Range("E2").Value = '  +1.90%  '
$ws.Range("D3").Value = '1.901.87'
$ws.Range("E3").Value = '  +3.07%  '
$ws.Range("E4").Value = '  +0.51%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.66'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.10%  '
$ws.Range("E6").Value = '  +2.57%  '
$ws.Range("E7").Value = '  +0.48%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '42.57'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.67%  '
$ws.Range("E9").Value = '  +2.91%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0707'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0992'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.77%  '
$ws.Range("D12").Value = '2.178.07'
$ws.Range("D14").Value = '1.892.30'
$ws.Range("E14").Value = '  +2.40%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.691'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.58%  '
$ws.Range("D17").Value = '35.616.18'
$ws.Range("E17").Value = '  +1.81%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '72.23'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.15%  '
$ws.Range("D19").Value = '0.0₃0811'
$ws.Range("E19").Value = '  +2.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '244.64'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.65%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.47'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.01%  '
$ws.Range("E22").Value = '  +3.58%  '
$ws.Range("E24").Value = '  +2.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '171.24'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.67%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.11'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +27.49%  '
$ws.Range("E27").Value = '  +7.90%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.76%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.958'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +28.28%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0567'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.54%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.09'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.24%  '
$ws.Range("E33").Value = '  +0.58%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.10'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.32%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.76'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.78%  '
$ws.Range("E36").Value = '  +4.38%  '
$ws.Range("E37").Value = '  +8.42%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.10'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.55%  '
$ws.Range("E39").Value = '  +4.91%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '91.44'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.90%  '
$ws.Range("D41").Value = '1.361.84'
$ws.Range("E41").Value = '  +0.88%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '15.39'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.50%  '
$ws.Range("E43").Value = '  +12.39%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.36'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.44%  '
$ws.Range("B45").Value = 'Gas'
$ws.Range("C45").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.10'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +36.11%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '47.17'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +39.10%  '
$ws.Range("E47").Value = '  +0.25%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.72'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.99%  '
$ws.Range("E49").Value = '  +0.59%  '
$ws.Range("D50").Value = '2.087.32'
$ws.Range("E50").Value = '  +2.69%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.52'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.44%  '
